# Auto-generated Excel COM-interop script
# Applies updated market-price / profit figures to the Anima_Profits leve tables
# across multiple crafting-class worksheets (ALC, ARM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets("ALC")
$ws.Range("H33").Value = 764.4167
$ws.Range("I33").Value = 1147.4286
$ws.Range("J33").Value = 228.2
$ws.Range("K33").Value = 1147.4286
$ws.Range("L33").Value = 228.2
$ws.Range("M33").Value = -918.4286
$ws.Range("N33").Value = -686.2

$ws.Range("H40").Value = 3277.7778
$ws.Range("I40").Value = 1500
$ws.Range("J40").Value = 3500
$ws.Range("K40").Value = 1500
$ws.Range("L40").Value = 3500
$ws.Range("M40").Value = -1325
$ws.Range("N40").Value = -3850

$ws.Range("H62").Value = 8917.272000000001
$ws.Range("I62").Value = 3460
$ws.Range("K62").Value = 3460
$ws.Range("M62").Value = -2836

$ws.Range("H64").Value = 2950.75
$ws.Range("J64").Value = 3001.2
$ws.Range("L64").Value = 3001.2
$ws.Range("N64").Value = -3497.2

$ws.Range("H65").Value = 8917.272000000001
$ws.Range("I65").Value = 3460
$ws.Range("K65").Value = 17300
$ws.Range("M65").Value = -14180

$ws.Range("H67").Value = 2950.75
$ws.Range("J67").Value = 3001.2
$ws.Range("L67").Value = 3001.2
$ws.Range("N67").Value = -4717.2

$ws.Range("H74").Value = 3990
$ws.Range("I74").Value = 3990
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 3990
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = $null
$ws.Range("N74").Value = -3054

$ws.Range("H77").Value = 3990
$ws.Range("I77").Value = 3990
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 19950
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = $null
$ws.Range("N77").Value = -15270

$ws.Range("H98").Value = 1211.625
$ws.Range("I98").Value = 1092.4375
$ws.Range("J98").Value = 1450
$ws.Range("K98").Value = 1092.4375
$ws.Range("L98").Value = 1450
$ws.Range("M98").Value = 405.5625
$ws.Range("N98").Value = -4446

$ws.Range("H122").Value = 1211.625
$ws.Range("I122").Value = 1092.4375
$ws.Range("J122").Value = 1450
$ws.Range("K122").Value = 3277.3125
$ws.Range("L122").Value = 4350
$ws.Range("M122").Value = -827.3125
$ws.Range("N122").Value = -9250

$ws.Range("H133").Value = 69780
$ws.Range("J133").Value = 69780
$ws.Range("L133").Value = 69780
$ws.Range("N133").Value = -79900

$ws.Range("H140").Value = 75715.75
$ws.Range("J140").Value = 75715.75
$ws.Range("L140").Value = 75715.75
$ws.Range("N140").Value = -86075.75

$ws = $wb.Sheets("ARM")
$ws.Range("H32").Value = 3947.68
$ws.Range("I32").Value = 3971.394
$ws.Range("J32").Value = 1600
$ws.Range("K32").Value = 3971.394
$ws.Range("L32").Value = 1600
$ws.Range("M32").Value = -3684.394
$ws.Range("N32").Value = -2174

$ws.Range("H122").Value = 73064.42999999999
$ws.Range("I122").Value = 78577.08
$ws.Range("J122").Value = 1400
$ws.Range("K122").Value = 235731.24
$ws.Range("L122").Value = 4200
$ws.Range("M122").Value = -233281.24
$ws.Range("N122").Value = -9100

$ws = $wb.Sheets("CRP")
$ws.Range("H31").Value = 6919.778
$ws.Range("I31").Value = 3262.75
$ws.Range("J31").Value = 8118.803
$ws.Range("K31").Value = 3262.75
$ws.Range("L31").Value = 8118.803
$ws.Range("M31").Value = -2967.75
$ws.Range("N31").Value = -8708.803

$ws.Range("H34").Value = 6919.778
$ws.Range("I34").Value = 3262.75
$ws.Range("J34").Value = 8118.803
$ws.Range("K34").Value = 3262.75
$ws.Range("L34").Value = 8118.803
$ws.Range("M34").Value = -3060.75
$ws.Range("N34").Value = -8522.803

$ws.Range("H93").Value = 3826.75
$ws.Range("I93").Value = 3826.75
$ws.Range("K93").Value = 3826.75
$ws.Range("M93").Value = -1954.75

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = $null
$ws.Range("M122").Value = $null
$ws.Range("N122").Value = 0

$ws = $wb.Sheets("CUL")
$ws.Range("H56").Value = 4837.778
$ws.Range("I56").Value = 4837.778
$ws.Range("K56").Value = 4837.778
$ws.Range("M56").Value = -4307.778

$ws.Range("H75").Value = 500
$ws.Range("I75").Value = 500
$ws.Range("K75").Value = 1500
$ws.Range("M75").Value = -502

$ws.Range("H78").Value = 500
$ws.Range("I78").Value = 500
$ws.Range("K78").Value = 4500
$ws.Range("M78").Value = 492

$ws.Range("H87").Value = 3260
$ws.Range("I87").Value = 2013.3334
$ws.Range("J87").Value = 7000
$ws.Range("K87").Value = 6040.0002
$ws.Range("L87").Value = 21000
$ws.Range("M87").Value = -4792.0002
$ws.Range("N87").Value = -23496

$ws.Range("H90").Value = 3260
$ws.Range("I90").Value = 2013.3334
$ws.Range("J90").Value = 7000
$ws.Range("K90").Value = 18120.0006
$ws.Range("L90").Value = 63000
$ws.Range("M90").Value = -11880.0006
$ws.Range("N90").Value = -75480

$ws.Range("H102").Value = 3490
$ws.Range("J102").Value = 3490
$ws.Range("L102").Value = 10470
$ws.Range("N102").Value = -15338

$ws.Range("H107").Value = 2408.3125
$ws.Range("I107").Value = 371.2857
$ws.Range("J107").Value = 3992.6667
$ws.Range("K107").Value = 1113.8571
$ws.Range("L107").Value = 11978.0001
$ws.Range("M107").Value = 806.1428999999998
$ws.Range("N107").Value = -15818.0001

$ws.Range("H121").Value = 803.1818
$ws.Range("I121").Value = 319.44446
$ws.Range("J121").Value = 2980
$ws.Range("K121").Value = 958.33338
$ws.Range("L121").Value = 8940
$ws.Range("M121").Value = 351.66662
$ws.Range("N121").Value = -11560

$ws.Range("H137").Value = 29638.684
$ws.Range("I137").Value = 5935.44
$ws.Range("J137").Value = 66675
$ws.Range("K137").Value = 17806.32
$ws.Range("L137").Value = 200025
$ws.Range("M137").Value = -12706.32
$ws.Range("N137").Value = -210225

$ws = $wb.Sheets("GSM")
$ws.Range("H122").Value = 3042.8572

$ws.Range("H133").Value = 60778
$ws.Range("J133").Value = 60778
$ws.Range("L133").Value = 60778
$ws.Range("N133").Value = -70898

$ws = $wb.Sheets("LTW")
$ws.Range("H38").Value = 333343330
$ws.Range("I38").Value = 333343330
$ws.Range("K38").Value = 333343330
$ws.Range("M38").Value = -333342920

$ws.Range("H40").Value = 3976.875
$ws.Range("I40").Value = 3973.5715
$ws.Range("K40").Value = 3973.5715
$ws.Range("M40").Value = -3837.5715

$ws.Range("H122").Value = 5530.326
$ws.Range("I122").Value = 4206.857
$ws.Range("J122").Value = 6109.3438
$ws.Range("K122").Value = 12620.571
$ws.Range("L122").Value = 18328.0314
$ws.Range("M122").Value = -10170.571
$ws.Range("N122").Value = -23228.0314

$ws = $wb.Sheets("WVR")
$ws.Range("H122").Value = 6000
$ws.Range("I122").Value = 4000
$ws.Range("J122").Value = 8000
$ws.Range("K122").Value = 12000
$ws.Range("L122").Value = 24000
$ws.Range("M122").Value = -9550
$ws.Range("N122").Value = -28900
